# Update the 'K' column (G) with the regenerated strikeout/K values.
# Per commit message: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals" -- the values in column G (header "K")
# are replaced with newly computed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 2
    4  = 2
    5  = 6
    6  = 8
    7  = 4
    8  = 0
    9  = 8
    10 = 11
    11 = 4
    12 = 4
    13 = 2
    14 = 3
    15 = 4
    16 = 0
    17 = 1
    18 = 2
    19 = 4
    20 = 2
    21 = 2
    22 = 5
    23 = 2
    24 = 2
    25 = 1
    26 = 0
    27 = 6
    28 = 1
    29 = 1
    30 = 1
    31 = 3
    32 = 3
    33 = 2
    34 = 4
    35 = 0
    36 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
